$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.680.65"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.899.83"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'578.81"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").Value = "'146.28"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "2.898.63"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'32.61"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "3.379.41"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").Value = "61.658.10"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "2.897.30"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'434.30"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'13.32"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "'0.659"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'6.95"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'79.95"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").Value = "'12.12"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "'10.22"
$ws.Range("E26").Value = "  -9.44%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'2.05"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("E29").Value = "  +17.92%  "
$ws.Range("D30").Value = "'7.13"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "'2.54"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'25.76"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("D36").Value = "'0.966"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("D37").Value = "'3.06"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").Value = "'5.49"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").Value = "'49.09"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").Value = "'1.97"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").Value = "'0.115"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("D44").Value = "'38.24"
$ws.Range("E44").Value = "  -4.20%  "
$ws.Range("D45").Value = "'134.80"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "2.684.57"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").Value = "'341.99"
$ws.Range("E48").Value = "  -5.78%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "'21.91"
$ws.Range("E51").Value = "  -4.82%  "
